# poisson_naive pronta para a rodada 27
#
# Renumbers the "Round" (col A) and "Matchweek" (col E) columns for every
# fixture row (2-27). Column E switches from a free-text "Matchweek N"
# label to a plain numeric N. Rows 23 and 24 additionally had their whole
# match record (columns B:BD) swapped between them (the two fixtures were
# re-ordered), on top of getting their own new A/E numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Round (A) / Matchweek (E) numbers, keyed by row.
$targets = @{
    2  = @(1, 1)
    3  = @(6, 4)
    4  = @(8, 6)
    5  = @(11, 8)
    6  = @(14, 10)
    7  = @(18, 13)
    8  = @(23, 17)
    9  = @(27, 19)
    10 = @(30, 21)
    11 = @(32, 22)
    12 = @(34, 15)
    13 = @(35, 24)
    14 = @(39, 26)
    15 = @(7, 5)
    16 = @(4, 3)
    17 = @(9, 7)
    18 = @(36, 23)
    19 = @(16, 12)
    20 = @(23, 14)
    21 = @(34, 25)
    22 = @(24, 18)
    23 = @(19, 20)
    24 = @(15, 11)
    25 = @(15, 16)
    26 = @(3, 2)
    27 = @(12, 9)
}

# Helper: write $value into $ref while preserving its "kind" (Excel
# auto-parses plain ISO date-looking strings like "2023-08-19" into date
# serials on assignment, which the source data never was - it's plain text).
function Set-CellLikeText($ws, $ref, $value) {
    if ($value -eq $null) {
        $value = ""
    }
    if ($value -is [string] -and $value -match '^\d{4}-\d{2}-\d{2}$') {
        # Force literal text via a quote prefix, then strip the resulting
        # quote-prefix style so the cell ends up as a plain text cell.
        $ws.Range($ref).Value = "'" + $value
        $ws.Range($ref).Style = "Normal"
    } else {
        $ws.Range($ref).Value = $value
    }
}

# --- Rows 23 & 24: swap their whole match record (columns B:BD) ---------
$dataCols = @("B","C","D","F","G","H","I","J","K","L","M","N","O","P","Q","R", `
              "S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF", `
              "AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS", `
              "AT","AU","AV","AW","AX","AY","AZ","BA","BB","BC","BD")

$row23vals = @{}
$row24vals = @{}
foreach ($c in $dataCols) {
    $row23vals[$c] = $ws.Range("$c`23").Value2
    $row24vals[$c] = $ws.Range("$c`24").Value2
}
foreach ($c in $dataCols) {
    Set-CellLikeText $ws "$c`23" $row24vals[$c]
    Set-CellLikeText $ws "$c`24" $row23vals[$c]
}

# --- All fixture rows: write the new Round (A) / Matchweek (E) numbers --
foreach ($row in ($targets.Keys | Sort-Object)) {
    $pair = $targets[$row]
    $ws.Range("A$row").Value = $pair[0]
    $ws.Range("E$row").Value = $pair[1]
}
